$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Covid-19 podatki")

# --- Carry the existing row formatting down onto the two new rows before ---
# --- writing values, so the new rows pick up the same look (fill/border/ ---
# --- number format) already used for the table's data rows.             ---
$ws.Range("A72:J72").Copy()
$ws.Range("A88:J88").PasteSpecial(-4122)
$ws.Range("A87:J87").Copy()
$ws.Range("A89:J89").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 88: 2020-06-06 ---
$ws.Range("A88").Value = 43988
$ws.Range("B88").Value = 83105
$ws.Range("C88").Value = 229
$ws.Range("D88").Value = 1485
$ws.Range("E88").Value = 1
$ws.Range("F88").Value = 5
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 1
$ws.Range("I88").Value = 109
$ws.Range("J88").Value = 0

# --- Row 89: 2020-06-07 ---
$ws.Range("A89").Value = 43989
$ws.Range("B89").Value = 83316
$ws.Range("C89").Value = 211
$ws.Range("D89").Value = 1485
$ws.Range("E89").Value = 0
$ws.Range("F89").Value = 5
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 109
$ws.Range("J89").Value = 0

# --- Grow the table / autofilter so the two new rows are part of "Tabela1" ---
$lo = $ws.ListObjects.Item("Tabela1")
$lo.Resize($ws.Range("A1:J89"))

# --- Match the saved view: last row selected, scrolled down ---
$null = $ws.Range("A89:J89").Select()
